$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 495; this shifts the existing rows 495:580 down to 496:581,
# preserving all of their data/styles untouched.
$ws.Rows.Item(495).Insert()

# Populate the newly inserted row 495 with its data.
$ws.Cells.Item(495, 1).Value = 5
$ws.Cells.Item(495, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(495, 3).Value = "Maule"
$ws.Cells.Item(495, 4).Value = 45244
$ws.Cells.Item(495, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(495, 5).Value = 7
$ws.Cells.Item(495, 6).Value = "Fruta"
$ws.Cells.Item(495, 7).Value = 100102
$ws.Cells.Item(495, 8).Value = "Cítricos"
$ws.Cells.Item(495, 9).Value = 100102004
$ws.Cells.Item(495, 10).Value = "Mandarina"
$ws.Cells.Item(495, 11).Value = "Murcott"
$ws.Cells.Item(495, 12).Value = "Primera"
$ws.Cells.Item(495, 13).Value = 300
$ws.Cells.Item(495, 14).Value = 8000
$ws.Cells.Item(495, 15).Value = 8000
$ws.Cells.Item(495, 16).Value = 8000
$ws.Cells.Item(495, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(495, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(495, 19).Value = 444
$ws.Cells.Item(495, 20).Value = 18
